$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column N (2020) mirroring the existing D:M year columns ---

# N3: same look as M3 (thin-bottom border cell, blank)
$ws.Range("M3").Copy($ws.Range("N3"))

# N4: same look as M4 (bold header, right/center, medium borders), new year value
$ws.Range("M4").Copy($ws.Range("N4"))
$ws.Range("N4").Value = 2020

# N5: looks similar to M5 (Times New Roman 9, medium top/bottom borders) but keeps
# General number format and only vertical-center alignment (no number formatting,
# no horizontal alignment) - matches a distinct new style.
$ws.Range("A5").Copy($ws.Range("N5"))
$ws.Range("N5").Value = 1.6
$ws.Range("N5").WrapText = $false
$ws.Range("N5").HorizontalAlignment = 1
$ws.Range("N5").NumberFormat = "General"

# --- Update existing value M5: 1.7 -> 1.6 ---
$ws.Range("M5").Value = 1.6

# --- Update current selection ---
$ws.Range("P6").Select()
